$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "2 - ..." step renumbered to "3 - ..." and reworded.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "2 - Usuário digita seu email, escolhe a categoria, descreve sua indagação e confirma ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "3 - Usuário digita seu email, escolhe a categoria, subcategoria e descreve sua indagação confirmando.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "7- ..." reworded ("sistema exibe mensagem" -> "é exibida uma mensagem").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "7- Usuário deixa campos sem preencher e confirma, então sistema exibe mensagem que os campos são obrigatórios e não envia mensagem;",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "7- Usuário deixa campos sem preencher e confirma, então é exibida uma mensagem que os campos são obrigatórios e não envia mensagem;",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "8- ..." sentence is reworded and split in two runs, with the
#    "_GoBack" bookmark moved from the document title ("COLLECTOR SHOP")
#    into the middle of this sentence (right after "uma").
# ---------------------------------------------------------------------------

# Find the split point: right after "...e exibe " in the original sentence.
$rFind = $d.Content
$rFind.Find.Execute("8- Sistema não consegue enviar mensagem e exibe ") | Out-Null
$leftStart = $rFind.Start
$leftEnd = $rFind.End

# Re-creating the bookmark under the same name ("_GoBack" is unique per
# document) removes the old one around "COLLECTOR SHOP" and places it here,
# splitting the run in two.
$bmRange = $d.Range($leftEnd, $leftEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Reword the left run (now its own run thanks to the bookmark split).
$rLeft = $d.Range($leftStart, $leftEnd)
$rLeft.Text = "8- Sistema não consegue enviar mensagem e exibe uma"

# Reword the right run.
$rRight = $d.Content
$rRight.Find.Execute("mensagem de erro. ") | Out-Null
$rRightRange = $d.Range($rRight.Start, $rRight.End)
$rRightRange.Text = " de erro. "

Write-Output "ok"
